$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("E2").Value = 25.72000000000058
$ws.Range("H2").Value = [double]"1.398706172756103e-16"
$ws.Range("K2").Value = 56.08575447588652
$ws.Range("L2").Value = "[52.14841146198581, 60.02309748978724]"
$ws.Range("O2").Value = 1.566079220708426
$ws.Range("P2").Value = "[1.490605523324887, 1.641552918091965]"
$ws.Range("S2").Value = 53.10600290122731
$ws.Range("T2").Value = "[50.389158904544615, 55.82284689791]"
$ws.Range("W2").Value = 19.30930930930974
$ws.Range("X2").Value = 19.00036036036078
$ws.Range("Y2").Value = 19.6182582582587

# Row 3 updates
$ws.Range("E3").Value = 24.92000000000046
$ws.Range("H3").Value = [double]"1.398706172756103e-16"
$ws.Range("K3").Value = 56.49438985732283
$ws.Range("L3").Value = "[49.55150808260954, 63.43727163203612]"
$ws.Range("O3").Value = 0.3710790121357315
$ws.Range("P3").Value = "[0.2327105669325782, 0.5094474573388847]"
$ws.Range("Q3").Value = [double]"4.328696017186218e-07"
$ws.Range("R3").Value = [double]"4.328696017186218e-07"
$ws.Range("S3").Value = 53.92042310189908
$ws.Range("T3").Value = "[49.99731210792006, 57.843534095878105]"
$ws.Range("W3").Value = 23.44824824824868
$ws.Range("X3").Value = 22.89945945945988
$ws.Range("Y3").Value = 23.99703703703747
